# Add a new "ExpectedFilenames" column (I) with expected report file-name
# prefixes for each study-type, used by the PROD test-data automation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "ExpectedFilenames"
$ws.Range("I2").Value = "Report-"
$ws.Range("I3").Value = "ExcelReport-NewImportLogic_1 - Test_Automation_1-Clinical-"
$ws.Range("I4").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Clinical-"
$ws.Range("I5").Value = "ExcelReport-NewImportLogic_1 - Test_Automation_1-Economic-"
$ws.Range("I6").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Economic-"
$ws.Range("I7").Value = "ExcelReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-"
$ws.Range("I8").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Quality of Life-"
$ws.Range("I9").Value = "ExcelReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-"
$ws.Range("I10").Value = "WordReport-NewImportLogic_1 - Test_Automation_1-Real-world Evidence-"

# The new column's cells should use the default (unstyled) format, matching
# the author's original edit, rather than inheriting the header/row styles.
$ws.Range("I1:I10").Style = "Normal"

# Leave the selection where the author finished editing.
$ws.Range("H10").Select()
